# "Generate Report for Handback" — refreshes the handback-status report:
#   - source record 34ba7391-71a2-4c5b-81f0-dc81687cc465 is renamed to 0d4c95d9-cb70-4025-8e6a-c91757c94803
#   - source record 444fc9ce-0507-4a2f-9cf6-b6ce109c4933 is renamed to ffff264d0ad0-d060-4efc-815b-4277b93913d1
#   - new handoff/handback timestamps and xliff checksums are recorded
#   - row 3 (ffff264d0ad0...) now shares its generated target file with row 2 (content duplicate)

$wb = $excel.ActiveWorkbook

# old ids (for reference): 34ba7391-71a2-4c5b-81f0-dc81687cc465, 444fc9ce-0507-4a2f-9cf6-b6ce109c4933
$newId1 = "0d4c95d9-cb70-4025-8e6a-c91757c94803"
$newId2 = "ffff264d0ad0-d060-4efc-815b-4277b93913d1"

$newChecksum = "af633a3421b149ad0631b37ac2d5e2870d361e85"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId1.md"
$wsOverview.Range("B2").Value = "e2e\$newId1.md"
$wsOverview.Range("G2").Value = "2016-08-21 05:07:38"

$wsOverview.Range("A3").Value = "$newId2.md"
$wsOverview.Range("B3").Value = "e2e\$newId2.md"
$wsOverview.Range("G3").Value = "2016-08-21 05:07:38"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newId1.md"
    }
    if ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$newId2.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId1.md"
$wsZhCn.Range("G2").Value = "$newId1.$newChecksum.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-21 05:07:34"
$wsZhCn.Range("I2").Value = "$newId1.md"
$wsZhCn.Range("J2").Value = "$newId1.$newChecksum.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-21 05:07:50"

$wsZhCn.Range("A3").Value = "$newId2.md"
$wsZhCn.Range("G3").Value = "$newId1.$newChecksum.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-21 05:07:34"
$wsZhCn.Range("I3").Value = "$newId2.md"
$wsZhCn.Range("J3").Value = "$newId1.$newChecksum.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-21 05:07:50"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newId1.md"
    }
    if ($addr -eq '$I$2') {
        $hl.TextToDisplay = "$newId1.md"
    }
    if ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newId2.md"
    }
    if ($addr -eq '$I$3') {
        $hl.TextToDisplay = "$newId2.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId1.md"
$wsDeDe.Range("G2").Value = "$newId1.$newChecksum.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-21 05:07:38"
$wsDeDe.Range("I2").Value = "$newId1.md"
$wsDeDe.Range("J2").Value = "$newId1.$newChecksum.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-21 05:07:57"

$wsDeDe.Range("A3").Value = "$newId2.md"
$wsDeDe.Range("G3").Value = "$newId1.$newChecksum.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-21 05:07:38"
$wsDeDe.Range("I3").Value = "$newId2.md"
$wsDeDe.Range("J3").Value = "$newId1.$newChecksum.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-21 05:07:57"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newId1.md"
    }
    if ($addr -eq '$I$2') {
        $hl.TextToDisplay = "$newId1.md"
    }
    if ($addr -eq '$A$3') {
        $hl.TextToDisplay = "$newId2.md"
    }
    if ($addr -eq '$I$3') {
        $hl.TextToDisplay = "$newId2.md"
    }
}
